$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '71.157.67'
$ws.Range("E2").Value = '  +0.30%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.810.59'
$ws.Range("E3").Value = '  -1.04%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '701.82'
$ws.Range("E5").Value = '  -0.59%  '
$ws.Range("E6").Value = '  -0.66%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.809.91'
$ws.Range("E7").Value = '  -1.01%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("E9").Value = '  -0.25%  '
$ws.Range("E10").Value = '  -0.77%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.50'
$ws.Range("E11").Value = '  +1.89%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.468'
$ws.Range("E12").Value = '  +1.97%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000252'
$ws.Range("E13").Value = '  -1.91%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.93'
$ws.Range("E14").Value = '  -1.93%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.451.89'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.846.43'
$ws.Range("E16").Value = '  +0.41%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '71.164.52'
$ws.Range("E17").Value = '  +0.21%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.51'
$ws.Range("E18").Value = '  +0.68%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.14'
$ws.Range("E19").Value = '  -0.80%  '
$ws.Range("E20").Value = '  -0.27%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '513.53'
$ws.Range("E21").Value = '  +4.17%  '
$ws.Range("E22").Value = '  -0.57%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.714'
$ws.Range("E23").Value = '  -0.42%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.92'
$ws.Range("E24").Value = '  -1.41%  '
$ws.Range("E25").Value = '  -3.22%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.958.28'
$ws.Range("E26").Value = '  -1.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.08'
$ws.Range("E27").Value = '  -0.65%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.37'
$ws.Range("E28").Value = '  -1.99%  '
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.04'
$ws.Range("E30").Value = '  -3.49%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.03'
$ws.Range("E31").Value = '  -4.73%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.40'
$ws.Range("E32").Value = '  -1.29%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.23'
$ws.Range("E33").Value = '  -1.56%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.02'
$ws.Range("E34").Value = '  -1.65%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.176'
$ws.Range("E35").Value = '  -2.56%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.12'
$ws.Range("E36").Value = '  -0.61%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.771.88'
$ws.Range("E37").Value = '  -0.91%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.990'
$ws.Range("E38").Value = '  -0.90%  '
$ws.Range("E39").Value = '  -2.62%  '
$ws.Range("E40").Value = '  +0.52%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.01'
$ws.Range("E41").Value = '  -1.07%  '
$ws.Range("E42").Value = '  -1.49%  '
$ws.Range("E43").Value = '  -1.24%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '172.87'
$ws.Range("E45").Value = '  +5.72%  '
$ws.Range("E46").Value = '  -0.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.000312'
$ws.Range("E47").Value = '  +1.11%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '49.35'
$ws.Range("E48").Value = '  +0.78%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '428.72'
$ws.Range("E49").Value = '  +3.39%  '
$ws.Range("E50").Value = '  -0.89%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.62'
$ws.Range("E51").Value = '  -0.16%  '
